$d = $word.ActiveDocument

# Locate the list item inside the instructions table that ends with
# "...with the appellate court clerk." (the first of the two occurrences
# of that phrase -- this one lives inside the table's last row/cell).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*File your form along with your other court papers with the appellate court clerk.*") {
        $target = $p
        break
    }
}

$full = $target.Range
$full.End = $full.End - 1   # exclude the paragraph mark
$full.Collapse(0)           # collapse to just after "clerk."

# Run 1: a plain space
$r1 = $full.Duplicate
$r1.InsertAfter(" ")
$r1.Font.Name = "Muli"
$r1.Font.Size = 14
$r1.Font.Bold = $false

# Run 2: bold "Note:"
$r2 = $d.Range($r1.End, $r1.End)
$r2.InsertAfter("Note:")
$r2.Font.Name = "Muli"
$r2.Font.Size = 14
$r2.Font.Bold = $true

# Run 3: the rest of the sentence
$r3 = $d.Range($r2.End, $r2.End)
$r3.InsertAfter(" Be sure to follow the filing requirements on the top of page 2 of the certification form.")
$r3.Font.Name = "Muli"
$r3.Font.Size = 14
$r3.Font.Bold = $false

Write-Output "Inserted Note sentence after target paragraph."
